{"js": "// Update the worksheet date and every \"three-digit \u00f7 one-digit\" answer cell\n// to the values generated for the new day (c986bee).\n//\n// Each (oldText -> newText) pair below is a UNIQUE paragraph/cell text run in\n// the source document, so a plain body-wide search-and-replace unambiguously\n// targets the correct run without needing to walk table coordinates.\nconst replacements = [\n  [\"2024-10-19 Saturday\", \"2024-10-20 Sunday\"],\n  [\"905\u00f79=100, 5\", \"641\u00f73=213, 2\"],\n  [\"830\u00f79=92, 2\", \"130\u00f73=43, 1\"],\n  [\"319\u00f74=79, 3\", \"101\u00f73=33, 2\"],\n  [\"355\u00f75=71, 0\", \"185\u00f78=23, 1\"],\n  [\"894\u00f76=149, 0\", \"522\u00f75=104, 2\"],\n  [\"505\u00f75=101, 0\", \"946\u00f78=118, 2\"],\n  [\"862\u00f75=172, 2\", \"927\u00f73=309, 0\"],\n  [\"660\u00f73=220, 0\", \"350\u00f79=38, 8\"],\n  [\"428\u00f79=47, 5\", \"700\u00f75=140, 0\"],\n  [\"360\u00f73=120, 0\", \"939\u00f72=469, 1\"],\n  [\"699\u00f78=87, 3\", \"660\u00f78=82, 4\"],\n  [\"633\u00f72=316, 1\", \"154\u00f74=38, 2\"],\n  [\"491\u00f78=61, 3\", \"639\u00f78=79, 7\"],\n  [\"815\u00f75=163, 0\", \"670\u00f79=74, 4\"],\n  [\"993\u00f77=141, 6\", \"778\u00f76=129, 4\"],\n  [\"245\u00f72=122, 1\", \"420\u00f77=60, 0\"],\n  [\"814\u00f79=90, 4\", \"807\u00f74=201, 3\"],\n  [\"406\u00f75=81, 1\", \"128\u00f73=42, 2\"],\n  [\"642\u00f76=107, 0\", \"375\u00f79=41, 6\"],\n  [\"697\u00f74=174, 1\", \"137\u00f78=17, 1\"],\n  [\"787\u00f75=157, 2\", \"209\u00f76=34, 5\"],\n  [\"819\u00f74=204, 3\", \"389\u00f76=64, 5\"],\n  [\"897\u00f73=299, 0\", \"283\u00f74=70, 3\"],\n  [\"894\u00f73=298, 0\", \"778\u00f75=155, 3\"],\n  [\"254\u00f74=63, 2\", \"652\u00f78=81, 4\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and every \"three-digit / one-digit\" answer cell\n# to the values generated for the new day (c986bee).\n#\n# Each (old -> new) pair is a unique text run in the source document, so a\n# body-wide Find/Replace (wdReplaceAll semantics via Replace:=2, but scoped to\n# a single exact hit each) unambiguously targets the correct run.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @('2024-10-19 Saturday', '2024-10-20 Sunday'),\n    @('905\u00f79=100, 5', '641\u00f73=213, 2'),\n    @('830\u00f79=92, 2', '130\u00f73=43, 1'),\n    @('319\u00f74=79, 3', '101\u00f73=33, 2'),\n    @('355\u00f75=71, 0', '185\u00f78=23, 1'),\n    @('894\u00f76=149, 0', '522\u00f75=104, 2'),\n    @('505\u00f75=101, 0', '946\u00f78=118, 2'),\n    @('862\u00f75=172, 2', '927\u00f73=309, 0'),\n    @('660\u00f73=220, 0', '350\u00f79=38, 8'),\n    @('428\u00f79=47, 5', '700\u00f75=140, 0'),\n    @('360\u00f73=120, 0', '939\u00f72=469, 1'),\n    @('699\u00f78=87, 3', '660\u00f78=82, 4'),\n    @('633\u00f72=316, 1', '154\u00f74=38, 2'),\n    @('491\u00f78=61, 3', '639\u00f78=79, 7'),\n    @('815\u00f75=163, 0', '670\u00f79=74, 4'),\n    @('993\u00f77=141, 6', '778\u00f76=129, 4'),\n    @('245\u00f72=122, 1', '420\u00f77=60, 0'),\n    @('814\u00f79=90, 4', '807\u00f74=201, 3'),\n    @('406\u00f75=81, 1', '128\u00f73=42, 2'),\n    @('642\u00f76=107, 0', '375\u00f79=41, 6'),\n    @('697\u00f74=174, 1', '137\u00f78=17, 1'),\n    @('787\u00f75=157, 2', '209\u00f76=34, 5'),\n    @('819\u00f74=204, 3', '389\u00f76=64, 5'),\n    @('897\u00f73=299, 0', '283\u00f74=70, 3'),\n    @('894\u00f73=298, 0', '778\u00f75=155, 3'),\n    @('254\u00f74=63, 2', '652\u00f78=81, 4')\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $true, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
